# Auto-generated Excel COM-interop script to apply the diff to before.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-21: only columns C-H (ax..gz) change; A (timestamp) / B (label) stay the same ---
$ws.Range("C2").Value = -2.236848592758179
$ws.Range("D2").Value = 3.817810773849488
$ws.Range("E2").Value = 1.356045484542847
$ws.Range("F2").Value = -0.0091629782691597
$ws.Range("G2").Value = 0.008399397134780801
$ws.Range("H2").Value = -0.0251981914043426

$ws.Range("C3").Value = -2.285426902770996
$ws.Range("D3").Value = 3.86536750793457
$ws.Range("E3").Value = 1.423394083976746
$ws.Range("F3").Value = 0.040775254368782
$ws.Range("G3").Value = 0.0474947728216648
$ws.Range("H3").Value = -0.0180205255746841

$ws.Range("C4").Value = -2.325414371490478
$ws.Range("D4").Value = 3.847443521022797
$ws.Range("E4").Value = 1.428170895576477
$ws.Range("F4").Value = 0.008399397134780801
$ws.Range("G4").Value = -0.0024434609804302
$ws.Range("H4").Value = -0.0125227374956011

$ws.Range("C5").Value = -2.198427677154541
$ws.Range("D5").Value = 3.85212025642395
$ws.Range("E5").Value = 1.377349805831909
$ws.Range("F5").Value = 0.0181732401251792
$ws.Range("G5").Value = 0.011148290708661
$ws.Range("H5").Value = 0.0201585534960031

$ws.Range("C6").Value = -2.239168739318848
$ws.Range("D6").Value = 3.77459921836853
$ws.Range("E6").Value = 1.236697590351105
$ws.Range("F6").Value = -0.0290160998702049
$ws.Range("G6").Value = -0.0050396383740007
$ws.Range("H6").Value = 0.0441350154578685

$ws.Range("C7").Value = -2.284214735031127
$ws.Range("D7").Value = 3.738507509231567
$ws.Range("E7").Value = 1.190666794776917
$ws.Range("F7").Value = 0.0397062413394451
$ws.Range("G7").Value = -0.0019853119738399
$ws.Range("H7").Value = 0.0597120784223079

$ws.Range("C8").Value = -2.274757814407349
$ws.Range("D8").Value = 3.780967509746552
$ws.Range("E8").Value = 1.303351855278015
$ws.Range("F8").Value = 0.0830776765942573
$ws.Range("G8").Value = -0.0916297882795333
$ws.Range("H8").Value = 0.06856962293386459

$ws.Range("C9").Value = -2.355447578430176
$ws.Range("D9").Value = 3.771934032440186
$ws.Range("E9").Value = 1.401894807815552
$ws.Range("F9").Value = 0.0453567430377006
$ws.Range("G9").Value = -0.0123700210824608
$ws.Range("H9").Value = 0.0334448739886283

$ws.Range("C10").Value = -2.428671193122864
$ws.Range("D10").Value = 3.703050696849823
$ws.Range("E10").Value = 1.323268264532089
$ws.Range("F10").Value = 0.0251981914043426
$ws.Range("G10").Value = 0.0175623763352632
$ws.Range("H10").Value = 0.0310014113783836

$ws.Range("C11").Value = -2.461624765396118
$ws.Range("D11").Value = 3.616514706611633
$ws.Range("E11").Value = 0.9642781019210801
$ws.Range("F11").Value = 0.0426078513264656
$ws.Range("G11").Value = -0.042302418500185
$ws.Range("H11").Value = 0.06704246252775189

$ws.Range("C12").Value = -2.319928169250488
$ws.Range("D12").Value = 3.585709452629088
$ws.Range("E12").Value = 0.5316097438335416
$ws.Range("F12").Value = -0.0807869285345077
$ws.Range("G12").Value = -0.06704246252775189
$ws.Range("H12").Value = -0.0733038261532783

$ws.Range("C13").Value = -2.305660724639893
$ws.Range("D13").Value = 3.332286834716796
$ws.Range("E13").Value = 0.07676682472228982
$ws.Range("F13").Value = -0.0183259565383195
$ws.Range("G13").Value = -0.1901318132877349
$ws.Range("H13").Value = 0.0618501044809818

$ws.Range("C14").Value = -2.763672423362733
$ws.Range("D14").Value = 2.987097477912901
$ws.Range("E14").Value = -0.2520411491394069
$ws.Range("F14").Value = -0.0526871271431446
$ws.Range("G14").Value = -0.3115412890911102
$ws.Range("H14").Value = 0.0861319974064827

$ws.Range("C15").Value = -3.228790092468267
$ws.Range("D15").Value = 2.132350564002988
$ws.Range("E15").Value = -1.067725944519043
$ws.Range("F15").Value = -0.057115901261568
$ws.Range("G15").Value = -0.2079996168613433
$ws.Range("H15").Value = -0.1476766765117645

$ws.Range("C16").Value = -5.195065975189209
$ws.Range("D16").Value = 1.12532408237457
$ws.Range("E16").Value = -1.079952371120452
$ws.Range("F16").Value = -0.0540615729987621
$ws.Range("G16").Value = 0.0456621758639812
$ws.Range("H16").Value = -0.3489567637443542

$ws.Range("C17").Value = -5.408200263977051
$ws.Range("D17").Value = 0.3617565631866455
$ws.Range("E17").Value = -0.6209573745727539
$ws.Range("F17").Value = 0.0517708286643028
$ws.Range("G17").Value = 0.1207986027002334
$ws.Range("H17").Value = -0.4889976382255554

$ws.Range("C18").Value = -5.044864058494567
$ws.Range("D18").Value = 0.9128529787063611
$ws.Range("E18").Value = -0.318287602066993
$ws.Range("F18").Value = 0.9677632451057434
$ws.Range("G18").Value = 1.55587375164032
$ws.Range("H18").Value = -2.213317394256592

$ws.Range("C19").Value = -2.147526121139514
$ws.Range("D19").Value = 2.947685408592233
$ws.Range("E19").Value = -1.418278175592428
$ws.Range("F19").Value = -0.3923282027244568
$ws.Range("G19").Value = 1.798387289047241
$ws.Range("H19").Value = -1.971720337867737

$ws.Range("C20").Value = -2.457507729530336
$ws.Range("D20").Value = 2.719761490821837
$ws.Range("E20").Value = -2.574102157354357
$ws.Range("F20").Value = -0.612850546836853
$ws.Range("G20").Value = -4.971374034881592
$ws.Range("H20").Value = -1.526399493217468

$ws.Range("C21").Value = -8.347633457183861
$ws.Range("D21").Value = -0.3747469902038705
$ws.Range("E21").Value = 2.810587501525903
$ws.Range("F21").Value = 0.1269072592258453
$ws.Range("G21").Value = 0.4891503453254699
$ws.Range("H21").Value = 0.1519527286291122

# --- Append 10 brand-new rows 22-31 (timestamps 2000-2900) ---
$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "falling"
$ws.Range("C22").Value = -9.06434655189514
$ws.Range("D22").Value = 2.370895385742192
$ws.Range("E22").Value = 2.526045709848401
$ws.Range("F22").Value = -0.0516181141138076
$ws.Range("G22").Value = 0.4092797040939331
$ws.Range("H22").Value = 0.1918116807937622

$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "falling"
$ws.Range("C23").Value = 1.765351390838669
$ws.Range("D23").Value = 3.193235373497008
$ws.Range("E23").Value = 0.3422038555145195
$ws.Range("F23").Value = -0.0490219369530677
$ws.Range("G23").Value = -0.1739438772201538
$ws.Range("H23").Value = 0.2171625941991806

$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "falling"
$ws.Range("C24").Value = 2.5670121669769
$ws.Range("D24").Value = 4.178837358951577
$ws.Range("E24").Value = 0.5715423285961221
$ws.Range("F24").Value = -0.0340557359158992
$ws.Range("G24").Value = -0.152105450630188
$ws.Range("H24").Value = -0.2504547536373138

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "falling"
$ws.Range("C25").Value = 0.4724056243896478
$ws.Range("D25").Value = 4.734174823760986
$ws.Range("E25").Value = 2.432305717468271
$ws.Range("F25").Value = -0.016951510682702
$ws.Range("G25").Value = -0.2273945808410644
$ws.Range("H25").Value = 0.3982841372489929

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "falling"
$ws.Range("C26").Value = 0.06286396980285353
$ws.Range("D26").Value = 4.10421558618545
$ws.Range("E26").Value = 2.57377957105636
$ws.Range("F26").Value = 0.0387899428606033
$ws.Range("G26").Value = -0.08124507963657369
$ws.Range("H26").Value = 0.1327104717493057

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "falling"
$ws.Range("C27").Value = -0.2147421836853027
$ws.Range("D27").Value = 3.490234732627868
$ws.Range("E27").Value = 1.606800019741057
$ws.Range("F27").Value = -0.0487165041267871
$ws.Range("G27").Value = -0.0951422601938247
$ws.Range("H27").Value = -0.1156062483787536

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "falling"
$ws.Range("C28").Value = 0.00004000663757422473
$ws.Range("D28").Value = 4.117116189002994
$ws.Range("E28").Value = 0.8414077341556527
$ws.Range("F28").Value = -0.0276416521519422
$ws.Range("G28").Value = 0.0578794814646244
$ws.Range("H28").Value = 0.0484110713005065

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "falling"
$ws.Range("C29").Value = 0.6697305679321306
$ws.Range("D29").Value = 4.689649581909178
$ws.Range("E29").Value = 0.7220227718353288
$ws.Range("F29").Value = 0.001527163083665
$ws.Range("G29").Value = 0.0532979927957057
$ws.Range("H29").Value = 0.061391957104206

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "falling"
$ws.Range("C30").Value = 1.128712320327757
$ws.Range("D30").Value = 4.191518974304199
$ws.Range("E30").Value = 1.193370014429091
$ws.Range("F30").Value = -0.008399397134780801
$ws.Range("G30").Value = -0.1032362282276153
$ws.Range("H30").Value = -0.090866208076477

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "falling"
$ws.Range("C31").Value = 0.9585402488708525
$ws.Range("D31").Value = 4.21038007736206
$ws.Range("E31").Value = 1.021172881126405
$ws.Range("F31").Value = 0.0235183127224445
$ws.Range("G31").Value = 0.0720821022987365
$ws.Range("H31").Value = -0.0161879286170005

